$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$values = @(
    "LA",
    "ENG ",
    "GER ",
    "FRE ",
    "SPA ",
    "ITA ",
    "RUS ",
    "CHI ",
    "UND ",
    "POR ",
    "POL ",
    "CZE ",
    "SLO ",
    "NOR ",
    "JPN ",
    "SWE ",
    "HRV ",
    "DUT ",
    "TUR ",
    "HUN ",
    "UKR ",
    "ROM ",
    "PERSIAN ",
    "GREC ",
    "FIN ",
    "DAN ",
    "ARABE "
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("A2:A27").Select()

$wb.Save()
